# Finalize BOM: fix "Total Cost per Board" formulas to multiply (price * qty)
# instead of add, add three new parts (TE terminal blocks + NXP mosfet), and
# add a bold/underlined "Total Cost per Board" summary row with the grand total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFmt = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"

# --- Fix existing rows 2-5: Total Cost per Board = Unit Price * Quantity ---
$ws.Range("E2").Formula = "=C2*D2"
$ws.Range("E3").Formula = "=C3*D3"
$ws.Range("E4").Formula = "=C4*D4"
$ws.Range("E5").Formula = "=C5*D5"

# --- Row 6: TE Connectivity 6P Terminal Block ---
$ws.Range("A6").Value = "TE   282837-6 "
$ws.Range("B6").Value = "TE Connectivity   6P Terminal Block"
$ws.Range("C6").Value = 3.18
$ws.Range("C6").NumberFormat = $currencyFmt
$ws.Range("D6").Value = 3
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("E6").NumberFormat = $currencyFmt
$ws.Range("F6").Value = "Digi-key"
$ws.Range("G6").Value = "https://www.digikey.com/product-detail/en/te-connectivity-amp-connectors/282837-6/A113322-ND/2187976"

# --- Row 7: TE Connectivity 2P Terminal Block ---
$ws.Range("A7").Value = "TE   282837-2"
$ws.Range("B7").Value = "TE Connectivity   2P Terminal Block"
$ws.Range("C7").Value = 1.04
$ws.Range("C7").NumberFormat = $currencyFmt
$ws.Range("D7").Value = 8
$ws.Range("E7").Formula = "=C7*D7"
$ws.Range("E7").NumberFormat = $currencyFmt
$ws.Range("F7").Value = "Digi-key"
$ws.Range("G7").Value = "https://www.digikey.com/product-detail/en/te-connectivity-amp-connectors/282837-2/A113320-ND/2187973"

# --- Row 8: NXP / Nexperia Mosfet ---
$ws.Range("A8").Value = "NXP PHP79NQ08LT,127"
$ws.Range("B8").Value = "Nexperia Mosfet"
$ws.Range("C8").Value = 1.08
$ws.Range("C8").NumberFormat = $currencyFmt
$ws.Range("D8").Value = 6
$ws.Range("E8").Formula = "=C8*D8"
$ws.Range("E8").NumberFormat = $currencyFmt
$ws.Range("F8").Value = "Digi-key"
$ws.Range("G8").Value = "https://www.digikey.com/products/en?keywords=PHP79NQ08LT"

# --- Row 10/11: Total Cost per Board summary ---
$ws.Range("E10").Value = "Total Cost per Board"
$ws.Range("E10").Font.Bold = $true
$ws.Range("E10").Font.Underline = $true

$ws.Range("E11").Formula = "=E2+E3+E4+E5+E6+E7+E8"
$ws.Range("E11").NumberFormat = $currencyFmt

# --- Cosmetic: widen columns A/B to fit new text (closest achievable widths) ---
$ws.Columns("A").ColumnWidth = 22.17
$ws.Columns("B").ColumnWidth = 30.0

# --- Cosmetic: move the selection like the author's last click before saving ---
[void]$ws.Range("B16").Select()
